$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '285.88'),
    @('E2', '0.70%'),
    @('D3', '29.36'),
    @('E3', '3.59%'),
    @('D4', '5.058'),
    @('E4', '0.46%'),
    @('D5', '0.06714'),
    @('E5', '3.16%'),
    @('D6', '7.318'),
    @('B7', 'GateToken'),
    @('C7', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D7', '3.442'),
    @('E7', '1.98%'),
    @('B8', 'FTXToken'),
    @('C8', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @('D8', '1.376'),
    @('E8', '-3.21%'),
    @('B9', 'MXToken'),
    @('C9', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D9', '0.9002'),
    @('E9', '-1.79%'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D10', '0.1598'),
    @('E10', '3.89%'),
    @('B11', 'LiechtensteinCryptoassetsExchange'),
    @('C11', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D11', '0.07057'),
    @('E11', '11.86%'),
    @('B12', 'MandalaExchangeToken'),
    @('C12', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D12', '0.07633'),
    @('E12', '2.01%'),
    @('B13', 'BitrueCoin'),
    @('C13', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D13', '0.02923'),
    @('E13', '3.04%'),
    @('B14', 'BitMartToken'),
    @('C14', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D14', '0.08984'),
    @('E14', '0.25%'),
    @('B15', 'BitForexToken'),
    @('C15', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D15', '0.001585'),
    @('E15', '-0.81%'),
    @('B16', 'CoinExToken'),
    @('C16', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
    @('D16', '0.04496'),
    @('E16', '1.40%'),
    @('B17', 'One'),
    @('C17', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('D17', '0.0006489'),
    @('E17', '0.92%'),
    @('B18', 'TigerCash'),
    @('C18', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D18', '0.006187'),
    @('E18', '1.58%'),
    @('B19', 'LEO'),
    @('C19', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D19', '3.447'),
    @('E19', '-0.02%'),
    @('D20', '2.232'),
    @('E20', '-0.30%'),
    @('E21', '1.30%'),
    @('D22', '0.1318'),
    @('E22', '0.82%'),
    @('D23', '3.908'),
    @('E23', '-2.21%'),
    @('E24', '1.63%'),
    @('D25', '0.001203'),
    @('E25', '1.63%'),
    @('D26', '0.004366'),
    @('E26', '-1.04%'),
    @('E27', '-6.48%'),
    @('E28', '-0.40%'),
    @('D40', '0.04238'),
    @('E40', '2.76%'),
    @('D41', '0.006769'),
    @('E41', '3.36%'),
    @('D42', '0.1237'),
    @('E42', '2.13%'),
    @('D43', '0.002229'),
    @('E43', '2.68%'),
    @('D44', '0.01262'),
    @('E44', '10.04%'),
    @('D45', '0.00005744'),
    @('E45', '0.94%'),
    @('D46', '1.963'),
    @('E46', '-0.14%'),
    @('E47', '15.07%')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
    $rng.Style = "Normal"
}
